$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted right before the
# existing row 267 (Membrillo / Primera / Región de O'Higgins), pushing
# every following row down by one. The new row repeats all the data from
# the (old) row 267 except for the date (column D) and the volume
# (column M), which carry the new week's figures.

$insertRow = 267

# Insert a blank row, shifting row 267 (and everything below it) down to 268.
$ws.Rows.Item($insertRow).Insert()

# The former row 267 now lives at row 268 - duplicate it back up into the
# freshly inserted row 267 so every other column matches.
$ws.Range("A268:T268").Copy()
$ws.Range("A267").PasteSpecial()

# Apply the two values that differ for the new entry: report date and volume.
$ws.Range("D267").Value = 45127
$ws.Range("M267").Value = 35
